$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Insert a new row at 124. This pushes the old "signature block" rows
#    (128 -> 129, 129 -> 130) down by one and extends the used range.
# ---------------------------------------------------------------------------
$ws.Rows.Item(124).Insert()

# ---------------------------------------------------------------------------
# 2) Row 123 used to be the LAST data row (special bottom-border style).
#    Now that row 124 exists, row 124 becomes the last row and should carry
#    that special style, while row 123 becomes a normal interior row.
# ---------------------------------------------------------------------------
$ws.Range("B123:J123").Copy()
$ws.Range("B124:J124").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B122:J122").Copy()
$ws.Range("B123:J123").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the values for the brand-new row 124 (same worker/amount columns
# as every other data row).
$ws.Range("B124").Value = "CC"
$ws.Range("C124").Value = "73114034"
$ws.Range("D124").Value = "GABRIEL ANTONIO HERRERA ARGUMEDO"
$ws.Range("F124").Value = 137360
$ws.Range("G124").Value = 3434000

# ---------------------------------------------------------------------------
# 3) Refresh the "Periodo Mora" column for every data row (16..124) so the
#    newest period (2507) is listed first and the data goes back further in
#    history (through 1607) instead of stopping at 2506.
# ---------------------------------------------------------------------------
$periods = @("2507","2506","2505","2504","2503","2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201","2112","2111","2110","2109","2108","2107","2106","2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801","1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701","1612","1611","1610","1609","1608","1607")
$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 5).Value = $p
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 4) Update the summary figures: total overdue amount and period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 14972240
$ws.Range("F13").Value = 109

Write-Host "Edit complete"
